# Log Week 17 data: update the "R" (road) row (row 3) on both the OFF and
# DEF sheets with the latest cumulative target-depth totals.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 239
$wsOff.Range("C3").Value = 166
$wsOff.Range("D3").Value = 59
$wsOff.Range("E3").Value = 27
$wsOff.Range("F3").Value = 3

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 151
$wsDef.Range("C3").Value = 103
$wsDef.Range("D3").Value = 36
$wsDef.Range("E3").Value = 17
